$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.139.43"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "1.841.06"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6256"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07487"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2939"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").Value = "1.878.83"
$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.020"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6757"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.975"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "29.154.90"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").Value = "2.128.86"
$ws.Range("E19").Value = "  +2.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.196"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.559"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1394"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.156"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05585"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.67%  "

$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7500"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.854"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.769"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").Value = "1.224.95"
$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01789"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.574"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").Value = "2.023.17"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5100"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4095"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.127"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05840"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.19%  "
